$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A261").Value = "IMX-USD"
$ws.Range("A262").Value = "MNT-USD"
$ws.Range("A263").Value = "TAO-USD"
